$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.051.87'
$ws.Range("E2").Value = '  -4.04%  '

$ws.Range("D3").Value = '1.962.60'
$ws.Range("E3").Value = '  -6.35%  '

$ws.Range("E4").Value = '  +0.50%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.97'

$ws.Range("E6").Value = '  +0.54%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4980'
$ws.Range("E7").Value = '  -5.64%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4208'
$ws.Range("E8").Value = '  -3.91%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.27'
$ws.Range("E9").Value = '  -1.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09026'
$ws.Range("E10").Value = '  -3.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.097'
$ws.Range("E11").Value = '  -6.47%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.99'
$ws.Range("E12").Value = '  -6.71%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.868'
$ws.Range("E13").Value = '  -8.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.439'
$ws.Range("E14").Value = '  -6.11%  '

$ws.Range("D15").Value = '1.923.12'
$ws.Range("E15").Value = '  -8.55%  '

$ws.Range("E16").Value = '  +0.66%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001097'
$ws.Range("E17").Value = '  -5.16%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.18'
$ws.Range("E18").Value = '  -9.70%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06648'
$ws.Range("E19").Value = '  -1.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.13'
$ws.Range("E20").Value = '  -9.49%  '

$ws.Range("E21").Value = '  +0.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.938'
$ws.Range("E22").Value = '  -6.94%  '

$ws.Range("D23").Value = '29.065.06'
$ws.Range("E23").Value = '  -4.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.89'
$ws.Range("E24").Value = '  -4.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.285'
$ws.Range("E25").Value = '  -1.52%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.58'
$ws.Range("E26").Value = '  -5.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.86'
$ws.Range("E27").Value = '  -4.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.146'
$ws.Range("E28").Value = '  -12.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.253'
$ws.Range("E29").Value = '  -10.19%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.20'
$ws.Range("E30").Value = '  -4.88%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.039'
$ws.Range("E31").Value = '  -8.25%  '

$ws.Range("E32").Value = '  -6.65%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.531'
$ws.Range("E33").Value = '  -8.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.788'
$ws.Range("E34").Value = '  -7.21%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.672'
$ws.Range("E35").Value = '  -6.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02420'
$ws.Range("E36").Value = '  -7.47%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.290'
$ws.Range("E37").Value = '  -3.63%  '

$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.914'
$ws.Range("E38").Value = '  -11.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06284'
$ws.Range("E39").Value = '  -6.99%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6439'
$ws.Range("E40").Value = '  -7.50%  '

$ws.Range("E41").Value = '  -9.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1984'
$ws.Range("E42").Value = '  -10.25%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.006'
$ws.Range("E43").Value = '  +0.46%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6186'
$ws.Range("E44").Value = '  -8.67%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.36'
$ws.Range("E45").Value = '  -6.92%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.164'
$ws.Range("E46").Value = '  -7.53%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.273'
$ws.Range("E47").Value = '  -2.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.461'
$ws.Range("E48").Value = '  -4.88%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000329'
$ws.Range("E49").Value = '  -6.31%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06885'
$ws.Range("E50").Value = '  -5.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.101'
$ws.Range("E51").Value = '  -8.84%  '
